$d = $word.ActiveDocument

# 1) Merge the two split runs describing the total de casos confirmados paragraph
$d.Content.Find.Execute(
    "O total de casos confirmados a que se refere este arquivo é a soma dos casos confirmados que não evoluíram para óbito e dos óbitos confirmados por COVID-19.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O total de casos confirmados a que se refere este arquivo é a soma dos casos confirmados que não evoluíram para óbito e dos óbitos confirmados por COVID-19.",
    2) | Out-Null

# 2) Merge the two split runs describing the "município de residência" paragraph
$d.Content.Find.Execute(
    "Este arquivo contém a relação dos casos confirmados com a especificação do sexo, idade, município de residência, ocorrência ou não de internações e se essas internações demandaram UTI.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Este arquivo contém a relação dos casos confirmados com a especificação do sexo, idade, município de residência, ocorrência ou não de internações e se essas internações demandaram UTI.",
    2) | Out-Null

# 3) Merge the two split runs describing the INTERNACAO_UTI field
$d.Content.Find.Execute(
    "Informa se o paciente com caso confirmado de COVID 19 internado, precisou de UTI (Unidade de Terapia Intensiva) com possibilidade de preenchimento SIM ou NÃO.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Informa se o paciente com caso confirmado de COVID 19 internado, precisou de UTI (Unidade de Terapia Intensiva) com possibilidade de preenchimento SIM ou NÃO.",
    2) | Out-Null

# 4) Table restructuring: row 8 (currently "7. | DATA_ATUALIZACAO | Data | Data do upload ...")
#    becomes the COMORBIDADE row, and a brand-new row is appended with the original
#    DATA_ATUALIZACAO content.
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(8)

$cell1 = $row.Cells.Item(1).Range
$cell2 = $row.Cells.Item(2).Range
$cell3 = $row.Cells.Item(3).Range
$cell4 = $row.Cells.Item(4).Range

# Append a brand-new row at the end of the table (copies formatting from row 8).
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "7."
$newRow.Cells.Item(2).Range.Text = "DATA_ATUALIZACAO"
$newRow.Cells.Item(3).Range.Text = "Data"
$newRow.Cells.Item(4).Range.Text = "Data do upload do arquivo, no formato YYYY-MM-DD."

# Now repurpose row 8 itself into the COMORBIDADE row.
$cell1.Text = "7"
$cell2.Text = "COMORBIDADE"
$cell3.Text = "Texto"
$cell4.Text = "Presença de doenças preexistentes/ comorbidades (diabetes, hipertensão, etc.), com possibilidade de preenchimento SIM ou NÃO."

# 5) Style update: mark "Fontepargpadro" (Default Paragraph Font) as semi-hidden.
$style = $d.Styles.Item("Fontepargpadro")
$style.SemiHidden = $true
